$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the computed / re-sampled values in rows 4-11 -----------------

$ws.Range("H4").Value2 = 8.42086361836311
$ws.Range("I4").Value2 = 0.501745677347789
$ws.Range("J4").Value2 = 30.3226130668644
$ws.Range("R4").Value2 = 6.4454799520623

$ws.Range("H5").Value2 = -17.3816690803508
$ws.Range("I5").Value2 = 0.551214899264904
$ws.Range("J5").Value2 = 15.7495765908952
$ws.Range("R5").Value2 = 3.03631159640739

$ws.Range("H6").Value2 = -28.8301741809225
$ws.Range("I6").Value2 = 0.599313443125452
$ws.Range("J6").Value2 = 1.18341586593766
$ws.Range("R6").Value2 = -2.52128815370404

$ws.Range("H7").Value2 = -38.1864536660425
$ws.Range("I7").Value2 = 0.650493343520364
$ws.Range("J7").Value2 = -9.65875567903714
$ws.Range("R7").Value2 = -13.38605712023

$ws.Range("H8").Value2 = -36.3227021182072
$ws.Range("I8").Value2 = 0.675648581248717
$ws.Range("J8").Value2 = -13.220752169069
$ws.Range("R8").Value2 = -18.4480822544681

$ws.Range("F9").Value2 = 46.6272780874236
$ws.Range("H9").Value2 = -9.55303290539078
$ws.Range("I9").Value2 = 0.700545785025913
$ws.Range("J9").Value2 = -13.5604847402453
$ws.Range("R9").Value2 = -20.201057990232

$ws.Range("F10").Value2 = 4.50871399467582
$ws.Range("H10").Value2 = 24.0593238332205
$ws.Range("I10").Value2 = 0.749644859282492
$ws.Range("J10").Value2 = -0.058513534306087
$ws.Range("R10").Value2 = -7.21322551882105

$ws.Range("I11").Value2 = 0.801917215853616
$ws.Range("J11").Value2 = 21.3194301272472

# --- Normalize the formatting of column R (now a "delta" column, like H/J) -
# Re-apply the (already effectively default) formatting properties so the
# cell gets re-bucketed into the default/base cell style, matching the
# H/J columns next to it.
$rCol = $ws.Range("R4:R10")
$rCol.HorizontalAlignment = 1
$rCol.VerticalAlignment = -4107
$rCol.WrapText = $false
$rCol.Orientation = 0
$rCol.IndentLevel = 0
$rCol.ShrinkToFit = $false
$rCol.Locked = $true
$rCol.FormulaHidden = $false

# --- Move/collapse the active selection from I3:J11 down to cell J11 -------
$ws.Range("J11").Select()
